# Add prepare functions for writing.
#
# Summary of the change (per commit "add prepare functions for writing"):
#  1. Rename the "rd_connect_test_dir" sheet/entity to
#     "rd_connect_test_directories" and its corresponding "entities" row
#     (name "dir" -> "directories").
#  2. Add a brand-new "deseases" entity (new row in "entities") backed by a
#     new, currently-empty sheet "rd_connect_test_deseases".
#  3. Point the existing attributes (date_of_inclusion, OrganizationID, type,
#     last_activities, name, ID) at the renamed "directories" entity.

$wb = $excel.ActiveWorkbook

$entitiesWs   = $wb.Worksheets.Item("entities")
$attributesWs = $wb.Worksheets.Item("attributes")
$dirWs        = $wb.Worksheets.Item("rd_connect_test_dir")

# 1. Rename the "dir" entity -> "directories" (sheet + entities row).
$dirWs.Name = "rd_connect_test_directories"
$entitiesWs.Range("A2").Value = "directories"

# 2. Append the new "deseases" entity row right below it.
$entitiesWs.Range("A3").Value = "deseases"
$entitiesWs.Range("B3").Value = "rd_connect_test"
$entitiesWs.Range("C3").Value = "Directory"
$entitiesWs.Range("G3").Value = "PostgreSQL"

# 3. Repoint the attributes at the renamed "directories" entity.
for ($r = 2; $r -le 7; $r++) {
    $attributesWs.Range("D" + $r).Value = "directories"
}

# 4. Add the new, still-empty "rd_connect_test_deseases" sheet as the last tab.
$newWs = $wb.Worksheets.Add($null, $wb.Worksheets($wb.Worksheets.Count))
$newWs.Name = "rd_connect_test_deseases"
